{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the renaming of ni_connect_nodes_to_nearest_point_on_nearest_edge_in_search\n// to ni_data_proc_connect_nodes_to_point_on_nearest_edge_in_search (title + SQL sample),\n// renames the Edge_table_prefix / Node_table_prefix parameters to\n// Edge_table_name / Node_table_name, and appends a new explanatory paragraph after the\n// figure caption.\n\nconst body = context.document.body;\n\nconst OLD_FN = \"ni_connect_nodes_to_nearest_point_on_nearest_edge_in_search\";\nconst NEW_FN = \"ni_data_proc_connect_nodes_to_point_on_nearest_edge_in_search\";\n\n// 1) Title occurrence (bold run at the very start of the document).\nconst titleResults = body.search(OLD_FN, { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(NEW_FN, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Parameter 1: \"Edge_table_prefix\" -> \"Edge_table_name\".\nconst edgeResults = body.search(\"Edge_table_prefix\", { matchCase: true });\nedgeResults.load(\"text\");\nawait context.sync();\nif (edgeResults.items.length > 0) {\n  edgeResults.items[0].insertText(\"Edge_table_name\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Parameter 4: \"Node_table_prefix\" -> \"Node_table_name\".\nconst nodeResults = body.search(\"Node_table_prefix\", { matchCase: true });\nnodeResults.load(\"text\");\nawait context.sync();\nif (nodeResults.items.length > 0) {\n  nodeResults.items[0].insertText(\"Node_table_name\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) SQL sample: rename the function call and add a space before the argument list.\nconst OLD_SQL = \"SELECT * FROM \" + OLD_FN + \"('data_national_grid_gas_pipeline_feeder','geom', '\";\nconst NEW_SQL = \"SELECT * FROM \" + NEW_FN + \" ('data_national_grid_gas_pipeline_feeder','geom', '\";\nconst sqlResults = body.search(OLD_SQL, { matchCase: true });\nsqlResults.load(\"text\");\nawait context.sync();\nif (sqlResults.items.length > 0) {\n  sqlResults.items[0].insertText(NEW_SQL, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 5) Append a new explanatory paragraph right after the Figure.1 caption paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet figurePara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Figure.1\") !== -1) {\n    figurePara = p;\n    break;\n  }\n}\nif (figurePara) {\n  const newPara = figurePara.insertParagraph(\n    \"This table would be output with _join appended to the input output table name (parameter 7). Secondly a table with _unique appended to the input output table name (parameter 7) is also written to the schema that contains the original geometry replaced with the newly derived geometry.\",\n    Word.InsertLocation.after\n  );\n  // The new paragraph should be plain text (not bold) like the rest of the body text.\n  newPara.font.bold = false;\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the renaming of ni_connect_nodes_to_nearest_point_on_nearest_edge_in_search\n# to ni_data_proc_connect_nodes_to_point_on_nearest_edge_in_search (title + SQL sample),\n# renames the Edge_table_prefix / Node_table_prefix parameters to\n# Edge_table_name / Node_table_name, and appends a new explanatory paragraph after the\n# figure caption.\n\n$d = $word.ActiveDocument\n\n$oldFn = \"ni_connect_nodes_to_nearest_point_on_nearest_edge_in_search\"\n$newFn = \"ni_data_proc_connect_nodes_to_point_on_nearest_edge_in_search\"\n\n# 1) Rename the function in the title (occurrence #1) and the SQL sample (occurrence #2),\n#    but leave the mention in the Figure.1 caption (occurrence #3) untouched.\n$renameRange = $d.Content\n$renameFind = $renameRange.Find\n$renameFind.ClearFormatting()\n$renameFind.Text = $oldFn\n$matchCount = 0\nwhile ($renameFind.Execute()) {\n  $matchCount = $matchCount + 1\n  if ($matchCount -le 2) {\n    $renameRange.Text = $newFn\n  }\n  $renameRange.Collapse(0)\n  $renameRange.SetRange($renameRange.End, $d.Content.End)\n}\n\n# 2) Add a space between the (now renamed) function name and the opening parenthesis in\n#    the SQL sample: \"...in_search('data...\" -> \"...in_search ('data...\".\n$parenRange = $d.Content\n$parenFind = $parenRange.Find\n$parenFind.ClearFormatting()\n$foundParen = $parenFind.Execute($newFn + \"(\")\nif ($foundParen) {\n  $parenRange.Text = $newFn + \" (\"\n}\n\n# 3) Parameter 1: \"Edge_table_prefix\" -> \"Edge_table_name\".\n$edgeRange = $d.Content\n$edgeFind = $edgeRange.Find\n$edgeFind.ClearFormatting()\n$foundEdge = $edgeFind.Execute(\"Edge_table_prefix\")\nif ($foundEdge) {\n  $edgeRange.Text = \"Edge_table_name\"\n}\n\n# 4) Parameter 4: \"Node_table_prefix\" -> \"Node_table_name\".\n$nodeRange = $d.Content\n$nodeFind = $nodeRange.Find\n$nodeFind.ClearFormatting()\n$foundNode = $nodeFind.Execute(\"Node_table_prefix\")\nif ($foundNode) {\n  $nodeRange.Text = \"Node_table_name\"\n}\n\n# 5) Append a new explanatory paragraph right after the Figure.1 caption paragraph.\n$paragraphs = $d.Paragraphs\n$figureIndex = -1\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n  if ($paragraphs.Item($i).Range.Text -like \"*Figure.1*\") {\n    $figureIndex = $i\n    break\n  }\n}\nif ($figureIndex -gt 0) {\n  $figurePara = $paragraphs.Item($figureIndex)\n  $figurePara.Range.InsertParagraphAfter()\n  $newPara = $d.Paragraphs.Item($figureIndex + 1)\n  $newPara.Range.Text = \"This table would be output with _join appended to the input output table name (parameter 7). Secondly a table with _unique appended to the input output table name (parameter 7) is also written to the schema that contains the original geometry replaced with the newly derived geometry.\"\n  # The new paragraph should be plain text (not bold) like the rest of the body text.\n  $newPara.Range.Font.Bold = 0\n}\n"}
